$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1: new header cell "Save" - copy the formatting from G1 (bold/border header style)
# then set the value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# H2, H3: new numeric data cells, default (unstyled) format.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
